# The "service line" for the Data Science team was renamed to "DS & BI".
# Update every player row on the "players" sheet that currently says
# "Data Science" so it reads "DS & BI" instead.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("players")

$used = $ws.UsedRange
for ($r = 1; $r -le $used.Rows.Count; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Text -eq "Data Science") {
        $cell.Value = "DS & BI"
    }
}

# Reflect the author's final UI state: the "players" tab ends up active,
# with cell C13 selected.
$ws.Activate()
$ws.Range("C13").Select()
